$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": refresh row 80 (precision fix on the timestamp, and make
# sure every cell carries its normal style) and append new row 81.
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Re-assert row 80's contents so the cells pick up explicit styling, and bump
# the timestamp's precision.
$wsAmsin.Rows.Item(80).ClearContents()

$wsAmsin.Cells.Item(80, 1).NumberFormat = "@"
$wsAmsin.Cells.Item(80, 1).Value = "2023-06-12"

$wsAmsin.Cells.Item(80, 2).Value = 45089.67525392361

$wsAmsin.Cells.Item(80, 3).Value = "178daytest"
$wsAmsin.Cells.Item(80, 4).Value = 124
$wsAmsin.Cells.Item(80, 5).Value = 124
$wsAmsin.Cells.Item(80, 6).Value = 0
$wsAmsin.Cells.Item(80, 7).Value = 1.9

# New row 81.
$wsAmsin.Cells.Item(81, 1).NumberFormat = "@"
$wsAmsin.Cells.Item(81, 1).Value = "2023-07-31"

$wsAmsin.Cells.Item(80, 2).Copy()
$wsAmsin.Cells.Item(81, 2).PasteSpecial(-4122)
$wsAmsin.Cells.Item(81, 2).Value = 45138.39590413195

$wsAmsin.Cells.Item(81, 3).Value = "180fnlrun"
$wsAmsin.Cells.Item(81, 4).Value = 124
$wsAmsin.Cells.Item(81, 5).Value = 118
$wsAmsin.Cells.Item(81, 6).Value = 6
$wsAmsin.Cells.Item(81, 7).Value = 2.23

# ---------------------------------------------------------------------------
# Sheet "BETA": append new row 37.
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Cells.Item(37, 1).NumberFormat = "@"
$wsBeta.Cells.Item(37, 1).Value = "2023-08-01"

$wsBeta.Cells.Item(36, 2).Copy()
$wsBeta.Cells.Item(37, 2).PasteSpecial(-4122)
$wsBeta.Cells.Item(37, 2).Value = 45139.54474763889

$wsBeta.Cells.Item(37, 3).Value = "180beta"
$wsBeta.Cells.Item(37, 4).Value = 124
$wsBeta.Cells.Item(37, 5).Value = 120
$wsBeta.Cells.Item(37, 6).Value = 4
$wsBeta.Cells.Item(37, 7).Value = 2.15

# ---------------------------------------------------------------------------
# Sheet "AMS": append new row 41.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Cells.Item(41, 1).NumberFormat = "@"
$wsAms.Cells.Item(41, 1).Value = "2023-08-01"
$wsAms.Cells.Item(39, 1).Copy()
$wsAms.Cells.Item(41, 1).PasteSpecial(-4122)

$wsAms.Cells.Item(40, 2).Copy()
$wsAms.Cells.Item(41, 2).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 2).Value = 45139.84310569906

$wsAms.Cells.Item(39, 3).Copy()
$wsAms.Cells.Item(41, 3).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 3).Value = "180live"

$wsAms.Cells.Item(39, 4).Copy()
$wsAms.Cells.Item(41, 4).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 4).Value = 124

$wsAms.Cells.Item(39, 5).Copy()
$wsAms.Cells.Item(41, 5).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 5).Value = 122

$wsAms.Cells.Item(39, 6).Copy()
$wsAms.Cells.Item(41, 6).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 6).Value = 2

$wsAms.Cells.Item(39, 7).Copy()
$wsAms.Cells.Item(41, 7).PasteSpecial(-4122)
$wsAms.Cells.Item(41, 7).Value = 1.81
